# Auto-generated edit script: updates market-price-derived columns (H-N)
# across all 8 sheets to match the refreshed values from the diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1659.1666
$ws.Range("I6").Value = 378.8889
$ws.Range("K6").Value = 1136.6667
$ws.Range("M6").Value = -1024.6667
$ws.Range("H18").Value = 35716868
$ws.Range("I18").Value = 45455904
$ws.Range("K18").Value = 45455904
$ws.Range("M18").Value = -45455620
$ws.Range("H28").Value = 47619412
$ws.Range("I28").Value = 47619412
$ws.Range("K28").Value = 47619412
$ws.Range("M28").Value = -47618927
$ws.Range("H51").Value = 3630.9355
$ws.Range("J51").Value = 3562.1538
$ws.Range("L51").Value = 3562.1538
$ws.Range("N51").Value = -4530.1538
$ws.Range("H57").Value = 135705.17
$ws.Range("J57").Value = 135705.17
$ws.Range("L57").Value = 407115.51
$ws.Range("N57").Value = -408113.51
$ws.Range("H69").Value = 30999.666
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("H72").Value = 30999.666
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H74").Value = 3600
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 3600
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("H86").Value = 12821.1
$ws.Range("I86").Value = 3085.125
$ws.Range("J86").Value = 51765
$ws.Range("K86").Value = 3085.125
$ws.Range("L86").Value = 51765
$ws.Range("M86").Value = -1962.125
$ws.Range("N86").Value = -54011
$ws.Range("H88").Value = 166716670
$ws.Range("I88").Value = 333333340
$ws.Range("J88").Value = 99999
$ws.Range("K88").Value = 333333340
$ws.Range("L88").Value = 99999
$ws.Range("M88").Value = -333332934
$ws.Range("N88").Value = -100811
$ws.Range("H89").Value = 12821.1
$ws.Range("I89").Value = 3085.125
$ws.Range("J89").Value = 51765
$ws.Range("K89").Value = 15425.625
$ws.Range("L89").Value = 258825
$ws.Range("M89").Value = -9809.625
$ws.Range("N89").Value = -270057
$ws.Range("H91").Value = 166716670
$ws.Range("I91").Value = 333333340
$ws.Range("J91").Value = 99999
$ws.Range("K91").Value = 333333340
$ws.Range("L91").Value = 99999
$ws.Range("M91").Value = -333331936
$ws.Range("N91").Value = -102807
$ws.Range("H93").Value = 25048.25
$ws.Range("J93").Value = 25048.25
$ws.Range("L93").Value = 25048.25
$ws.Range("N93").Value = -30040.25
$ws.Range("H95").Value = 24415
$ws.Range("J95").Value = 24415
$ws.Range("L95").Value = 24415
$ws.Range("N95").Value = -29907
$ws.Range("H98").Value = 751.125
$ws.Range("I98").Value = 805.6923
$ws.Range("K98").Value = 805.6923
$ws.Range("M98").Value = 692.3077
$ws.Range("H106").Value = 19610068
$ws.Range("I106").Value = 22224518
$ws.Range("K106").Value = 22224518
$ws.Range("M106").Value = -22223887
$ws.Range("H107").Value = 20844720
$ws.Range("I107").Value = 21750142
$ws.Range("K107").Value = 21750142
$ws.Range("M107").Value = -21748222
$ws.Range("H113").Value = 1001
$ws.Range("I113").Value = 1001
$ws.Range("K113").Value = 1001
$ws.Range("M113").Value = 2253
$ws.Range("H116").Value = 9124.875
$ws.Range("I116").Value = 32999
$ws.Range("J116").Value = 5714.2856
$ws.Range("K116").Value = 32999
$ws.Range("L116").Value = 5714.2856
$ws.Range("M116").Value = -29557
$ws.Range("N116").Value = -12598.2856
$ws.Range("H122").Value = 751.125
$ws.Range("I122").Value = 805.6923
$ws.Range("K122").Value = 2417.0769
$ws.Range("M122").Value = 32.92309999999998
$ws.Range("H131").Value = 6222.9395
$ws.Range("I131").Value = 1675.1428
$ws.Range("J131").Value = 9573.947
$ws.Range("K131").Value = 5025.428400000001
$ws.Range("L131").Value = 28721.841
$ws.Range("M131").Value = 14.57159999999931
$ws.Range("N131").Value = -38801.841
$ws.Range("H136").Value = 179998
$ws.Range("J136").Value = 179998
$ws.Range("L136").Value = 179998
$ws.Range("N136").Value = -190198
$ws.Range("H137").Value = 1515.8572
$ws.Range("I137").Value = 1101.8334
$ws.Range("K137").Value = 3305.5002
$ws.Range("M137").Value = -755.5001999999999
$ws.Range("H138").Value = 1526.9474
$ws.Range("I138").Value = 1237.1212
$ws.Range("J138").Value = 3439.8
$ws.Range("K138").Value = 3711.3636
$ws.Range("L138").Value = 10319.4
$ws.Range("M138").Value = 1428.6364
$ws.Range("N138").Value = -20599.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 13911542
$ws.Range("I2").Value = 15547939
$ws.Range("J2").Value = 2171.5
$ws.Range("K2").Value = 15547939
$ws.Range("L2").Value = 2171.5
$ws.Range("M2").Value = -15547826
$ws.Range("N2").Value = -2397.5
$ws.Range("H32").Value = 4396.4287
$ws.Range("I32").Value = 3156.3103
$ws.Range("K32").Value = 3156.3103
$ws.Range("M32").Value = -2869.3103
$ws.Range("H45").Value = 19899.572
$ws.Range("I45").Value = 19899.572
$ws.Range("K45").Value = 19899.572
$ws.Range("M45").Value = -19522.572
$ws.Range("H61").Value = 9105.825999999999
$ws.Range("I61").Value = 8221.75
$ws.Range("K61").Value = 8221.75
$ws.Range("M61").Value = -8009.75
$ws.Range("H74").Value = 5668.8647
$ws.Range("I74").Value = 5867.5938
$ws.Range("K74").Value = 5867.5938
$ws.Range("M74").Value = -4993.5938
$ws.Range("H77").Value = 5668.8647
$ws.Range("I77").Value = 5867.5938
$ws.Range("K77").Value = 29337.969
$ws.Range("M77").Value = -24969.969
$ws.Range("H97").Value = 35720556
$ws.Range("I97").Value = 47625330
$ws.Range("K97").Value = 47625330
$ws.Range("M97").Value = -47624834
$ws.Range("H116").Value = 13911542
$ws.Range("I116").Value = 15547939
$ws.Range("J116").Value = 2171.5
$ws.Range("K116").Value = 15547939
$ws.Range("L116").Value = 2171.5
$ws.Range("M116").Value = -15545645
$ws.Range("N116").Value = -6759.5
$ws.Range("H122").Value = 4545
$ws.Range("I122").Value = 3660.8333
$ws.Range("K122").Value = 10982.4999
$ws.Range("M122").Value = -8532.499899999999
$ws.Range("H136").Value = 9105.825999999999
$ws.Range("I136").Value = 8221.75
$ws.Range("K136").Value = 24665.25
$ws.Range("M136").Value = -22115.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 13911542
$ws.Range("I3").Value = 15547939
$ws.Range("J3").Value = 2171.5
$ws.Range("K3").Value = 15547939
$ws.Range("L3").Value = 2171.5
$ws.Range("M3").Value = -15547825
$ws.Range("N3").Value = -2399.5
$ws.Range("H20").Value = 1638.3334
$ws.Range("I20").Value = 1256.9143
$ws.Range("J20").Value = 2665.2307
$ws.Range("K20").Value = 1256.9143
$ws.Range("L20").Value = 2665.2307
$ws.Range("M20").Value = -1009.9143
$ws.Range("N20").Value = -3159.2307
$ws.Range("H25").Value = 12500
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("H108").Value = 49999.9
$ws.Range("J108").Value = 49999.9
$ws.Range("L108").Value = 49999.9
$ws.Range("N108").Value = -57679.9
$ws.Range("H132").Value = 101304
$ws.Range("J132").Value = 101304
$ws.Range("L132").Value = 101304
$ws.Range("N132").Value = -111424
$ws.Range("H134").Value = 4552.2095
$ws.Range("I134").Value = 4503.9736
$ws.Range("J134").Value = 4918.8
$ws.Range("K134").Value = 13511.9208
$ws.Range("L134").Value = 14756.4
$ws.Range("M134").Value = -10976.9208
$ws.Range("N134").Value = -19826.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 48.857143
$ws.Range("I7").Value = 22.5
$ws.Range("J7").Value = 84
$ws.Range("K7").Value = 22.5
$ws.Range("L7").Value = 84
$ws.Range("M7").Value = 90.5
$ws.Range("N7").Value = -310
$ws.Range("H16").Value = 4778.8
$ws.Range("J16").Value = 4666.6665
$ws.Range("L16").Value = 4666.6665
$ws.Range("N16").Value = -5240.6665
$ws.Range("H22").Value = 2130.4666
$ws.Range("I22").Value = 2317
$ws.Range("K22").Value = 2317
$ws.Range("M22").Value = -1967
$ws.Range("H58").Value = 7240.517
$ws.Range("I58").Value = 8399.691999999999
$ws.Range("J58").Value = 6298.6875
$ws.Range("K58").Value = 8399.691999999999
$ws.Range("L58").Value = 6298.6875
$ws.Range("M58").Value = -8196.691999999999
$ws.Range("N58").Value = -6704.6875
$ws.Range("H105").Value = 2425.6428
$ws.Range("I105").Value = 2414.5454
$ws.Range("J105").Value = 2466.3333
$ws.Range("K105").Value = 2414.5454
$ws.Range("L105").Value = 2466.3333
$ws.Range("M105").Value = -667.5454
$ws.Range("N105").Value = -5960.3333
$ws.Range("H113").Value = 4778.8
$ws.Range("J113").Value = 4666.6665
$ws.Range("L113").Value = 4666.6665
$ws.Range("N113").Value = -9006.666499999999
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H122").Value = 3768.4546
$ws.Range("I122").Value = 2798.8333
$ws.Range("J122").Value = 4932
$ws.Range("K122").Value = 8396.499899999999
$ws.Range("L122").Value = 14796
$ws.Range("M122").Value = -5946.499899999999
$ws.Range("N122").Value = -19696
$ws.Range("H132").Value = 5002.173
$ws.Range("I132").Value = 5540.024
$ws.Range("J132").Value = 2743.2
$ws.Range("K132").Value = 16620.072
$ws.Range("L132").Value = 8229.599999999999
$ws.Range("M132").Value = -14090.072
$ws.Range("N132").Value = -13289.6
$ws.Range("H136").Value = 7240.517
$ws.Range("I136").Value = 8399.691999999999
$ws.Range("J136").Value = 6298.6875
$ws.Range("K136").Value = 25199.076
$ws.Range("L136").Value = 18896.0625
$ws.Range("M136").Value = -22649.076
$ws.Range("N136").Value = -23996.0625
$ws.Range("H138").Value = 123240.125
$ws.Range("J138").Value = 123240.125
$ws.Range("L138").Value = 123240.125
$ws.Range("N138").Value = -133520.125
$ws.Range("H140").Value = 108107.336
$ws.Range("I140").Value = 104790.75
$ws.Range("J140").Value = 114740.5
$ws.Range("K140").Value = 104790.75
$ws.Range("L140").Value = 114740.5
$ws.Range("M140").Value = -99610.75
$ws.Range("N140").Value = -125100.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3092.077
$ws.Range("I3").Value = 1836.091
$ws.Range("K3").Value = 5508.272999999999
$ws.Range("M3").Value = -5396.272999999999
$ws.Range("H12").Value = 49.333332
$ws.Range("J12").Value = 48.18182
$ws.Range("L12").Value = 144.54546
$ws.Range("N12").Value = -490.54546
$ws.Range("H33").Value = 392.27274
$ws.Range("I33").Value = 199.66667
$ws.Range("K33").Value = 1198.00002
$ws.Range("M33").Value = -915.0000199999999
$ws.Range("H58").Value = 222.5
$ws.Range("I58").Value = 222.5
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 667.5
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -539.5
$ws.Range("N58").ClearContents()
$ws.Range("H132").Value = 2461.2354
$ws.Range("I132").Value = 2195.6365
$ws.Range("J132").Value = 2948.1667
$ws.Range("K132").Value = 19760.7285
$ws.Range("L132").Value = 26533.5003
$ws.Range("M132").Value = -17230.7285
$ws.Range("N132").Value = -31593.5003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 24283.5
$ws.Range("I41").Value = 27540.2
$ws.Range("K41").Value = 27540.2
$ws.Range("M41").Value = -27185.2
$ws.Range("H70").Value = 4119.9395
$ws.Range("I70").Value = 3848.4119
$ws.Range("J70").Value = 4408.4375
$ws.Range("K70").Value = 3848.4119
$ws.Range("L70").Value = 4408.4375
$ws.Range("M70").Value = -3578.4119
$ws.Range("N70").Value = -4948.4375
$ws.Range("H73").Value = 4119.9395
$ws.Range("I73").Value = 3848.4119
$ws.Range("J73").Value = 4408.4375
$ws.Range("K73").Value = 3848.4119
$ws.Range("L73").Value = 4408.4375
$ws.Range("M73").Value = -2912.4119
$ws.Range("N73").Value = -6280.4375
$ws.Range("H80").Value = 34290196
$ws.Range("I80").Value = 60003344
$ws.Range("J80").Value = 5998.6665
$ws.Range("K80").Value = 60003344
$ws.Range("L80").Value = 5998.6665
$ws.Range("M80").Value = -60002346
$ws.Range("N80").Value = -7994.6665
$ws.Range("H83").Value = 34290196
$ws.Range("I83").Value = 60003344
$ws.Range("J83").Value = 5998.6665
$ws.Range("K83").Value = 300016720
$ws.Range("L83").Value = 29993.3325
$ws.Range("M83").Value = -300011728
$ws.Range("N83").Value = -39977.3325
$ws.Range("H97").Value = 1185.1666
$ws.Range("I97").Value = 886.6667
$ws.Range("K97").Value = 886.6667
$ws.Range("M97").Value = -390.6667
$ws.Range("H113").Value = 7580687
$ws.Range("I113").Value = 5854.5
$ws.Range("J113").Value = 9263983
$ws.Range("K113").Value = 5854.5
$ws.Range("L113").Value = 9263983
$ws.Range("M113").Value = -3684.5
$ws.Range("N113").Value = -9268323
$ws.Range("H122").Value = 4576.5
$ws.Range("I122").Value = 4881
$ws.Range("J122").Value = 4069
$ws.Range("K122").Value = 14643
$ws.Range("L122").Value = 12207
$ws.Range("M122").Value = -12193
$ws.Range("N122").Value = -17107
$ws.Range("H132").Value = 3808.8667
$ws.Range("I132").Value = 4053.7585
$ws.Range("J132").Value = 3365
$ws.Range("K132").Value = 12161.2755
$ws.Range("L132").Value = 10095
$ws.Range("M132").Value = -9631.2755
$ws.Range("N132").Value = -15155
$ws.Range("H135").Value = 133485
$ws.Range("J135").Value = 139066
$ws.Range("L135").Value = 139066
$ws.Range("N135").Value = -149206

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2651.3333
$ws.Range("I7").Value = 2636.2307
$ws.Range("K7").Value = 2636.2307
$ws.Range("M7").Value = -2524.2307
$ws.Range("H31").Value = 4273.4
$ws.Range("I31").Value = 5900
$ws.Range("J31").Value = 3576.2856
$ws.Range("K31").Value = 5900
$ws.Range("L31").Value = 3576.2856
$ws.Range("M31").Value = -5652
$ws.Range("N31").Value = -4072.2856
$ws.Range("H46").Value = 2969.5217
$ws.Range("I46").Value = 1846.4546
$ws.Range("J46").Value = 3999
$ws.Range("K46").Value = 1846.4546
$ws.Range("L46").Value = 3999
$ws.Range("M46").Value = -1658.4546
$ws.Range("N46").Value = -4375
$ws.Range("H93").Value = 25000792
$ws.Range("I93").Value = 28571978
$ws.Range("K93").Value = 28571978
$ws.Range("M93").Value = -28570730
$ws.Range("H115").Value = 149998
$ws.Range("J115").Value = 149998
$ws.Range("L115").Value = 149998
$ws.Range("N115").Value = -152348
$ws.Range("H122").Value = 11531.923
$ws.Range("I122").Value = 12567.889
$ws.Range("K122").Value = 37703.667
$ws.Range("M122").Value = -35253.667
$ws.Range("H126").Value = 2651.3333
$ws.Range("I126").Value = 2636.2307
$ws.Range("K126").Value = 7908.6921
$ws.Range("M126").Value = -5438.6921
$ws.Range("H132").Value = 25669.432
$ws.Range("I132").Value = 28206.111
$ws.Range("J132").Value = 6644.3335
$ws.Range("K132").Value = 84618.333
$ws.Range("L132").Value = 19933.0005
$ws.Range("M132").Value = -82088.333
$ws.Range("N132").Value = -24993.0005
$ws.Range("H136").Value = 4508061.5
$ws.Range("I136").Value = 5460453.5
$ws.Range("J136").Value = 18213.428
$ws.Range("K136").Value = 16381360.5
$ws.Range("L136").Value = 54640.284
$ws.Range("M136").Value = -16378810.5
$ws.Range("N136").Value = -59740.284
$ws.Range("H138").Value = 104959.43
$ws.Range("J138").Value = 104959.43
$ws.Range("L138").Value = 104959.43
$ws.Range("N138").Value = -115239.43
$ws.Range("H141").Value = 140710.67
$ws.Range("J141").Value = 140710.67
$ws.Range("L141").Value = 140710.67
$ws.Range("N141").Value = -151070.67

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 257599.8
$ws.Range("J15").Value = 129000
$ws.Range("L15").Value = 129000
$ws.Range("N15").Value = -129576
$ws.Range("H61").Value = 5000
$ws.Range("I61").Value = 5000
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 5000
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -4708
$ws.Range("N61").ClearContents()
$ws.Range("H62").Value = 76926910
$ws.Range("I62").Value = 125003270
$ws.Range("K62").Value = 125003270
$ws.Range("M62").Value = -125002646
$ws.Range("H65").Value = 76926910
$ws.Range("I65").Value = 125003270
$ws.Range("K65").Value = 625016350
$ws.Range("M65").Value = -625013230
$ws.Range("H81").Value = 4792064.5
$ws.Range("I81").Value = 7579102
$ws.Range("J81").Value = 14285.571
$ws.Range("K81").Value = 15158204
$ws.Range("L81").Value = 28571.142
$ws.Range("M81").Value = -15157143
$ws.Range("N81").Value = -30693.142
$ws.Range("H84").Value = 4792064.5
$ws.Range("I84").Value = 7579102
$ws.Range("J84").Value = 14285.571
$ws.Range("K84").Value = 75791020
$ws.Range("L84").Value = 142855.71
$ws.Range("M84").Value = -75785716
$ws.Range("N84").Value = -153463.71
$ws.Range("H96").Value = 3618.0833
$ws.Range("I96").Value = 4878
$ws.Range("J96").Value = 2358.1667
$ws.Range("K96").Value = 4878
$ws.Range("L96").Value = 2358.1667
$ws.Range("M96").Value = -3505
$ws.Range("N96").Value = -5104.1667
$ws.Range("H100").Value = 696.3125
$ws.Range("I100").Value = 595.5454999999999
$ws.Range("K100").Value = 1191.091
$ws.Range("M100").Value = -650.0909999999999
$ws.Range("H112").Value = 89999.664
$ws.Range("J112").Value = 89999.664
$ws.Range("L112").Value = 89999.664
$ws.Range("N112").Value = -92953.664
$ws.Range("H115").Value = 35999.25
$ws.Range("J115").Value = 35999.25
$ws.Range("L115").Value = 35999.25
$ws.Range("N115").Value = -39133.25
$ws.Range("H122").Value = 5000.926
$ws.Range("I122").Value = 4524.423
$ws.Range("J122").Value = 5443.393
$ws.Range("K122").Value = 13573.269
$ws.Range("L122").Value = 16330.179
$ws.Range("M122").Value = -11123.269
$ws.Range("N122").Value = -21230.179
$ws.Range("H126").Value = 6865.9165
$ws.Range("I126").Value = 6533.8335
$ws.Range("J126").Value = 7530.0835
$ws.Range("K126").Value = 19601.5005
$ws.Range("L126").Value = 22590.2505
$ws.Range("M126").Value = -17131.5005
$ws.Range("N126").Value = -27530.2505
$ws.Range("H132").Value = 1646.8667
$ws.Range("I132").Value = 1373
$ws.Range("J132").Value = 2742.3333
$ws.Range("K132").Value = 4119
$ws.Range("L132").Value = 8226.999899999999
$ws.Range("M132").Value = -1589
$ws.Range("N132").Value = -13286.9999
